$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'FC Bayern München'
$ws.Cells.Item(2, 2).Value = 5.691950464396285
$ws.Cells.Item(2, 3).Value = 8.486024844720497
$ws.Cells.Item(2, 4).Value = 0.806721111424012
$ws.Cells.Item(2, 5).Value = 24
$ws.Cells.Item(2, 6).Value = 0.04927499037597844
$ws.Cells.Item(2, 7).Value = 75

$ws.Cells.Item(3, 1).Value = 'Borussia Dortmund'
$ws.Cells.Item(3, 2).Value = 5.840755735492578
$ws.Cells.Item(3, 3).Value = 7.206666666666667
$ws.Cells.Item(3, 4).Value = 0.6086585365853658
$ws.Cells.Item(3, 5).Value = 14
$ws.Cells.Item(3, 6).Value = 0.07163557326571274
$ws.Cells.Item(3, 7).Value = 20

$ws.Cells.Item(4, 1).Value = 'TSG Hoffenheim'
$ws.Cells.Item(4, 2).Value = 4.561382598331347
$ws.Cells.Item(4, 3).Value = 8.403883495145632
$ws.Cells.Item(4, 4).Value = 0.5543060498220641
$ws.Cells.Item(4, 5).Value = 11
$ws.Cells.Item(4, 6).Value = 0.08453184619900206
$ws.Cells.Item(4, 7).Value = 12

$ws.Cells.Item(5, 1).Value = 'VfB Stuttgart'
$ws.Cells.Item(5, 2).Value = 4.853424657534247
$ws.Cells.Item(5, 3).Value = 8.082397003745319
$ws.Cells.Item(5, 4).Value = 0.6492954324586978
$ws.Cells.Item(5, 5).Value = 13
$ws.Cells.Item(5, 6).Value = 0.05803063457330416
$ws.Cells.Item(5, 7).Value = 13

$ws.Cells.Item(6, 1).Value = 'RB Leipzig'
$ws.Cells.Item(6, 2).Value = 5.807639836289223
$ws.Cells.Item(6, 3).Value = 7.919129082426127
$ws.Cells.Item(6, 4).Value = 0.5463829787234042
$ws.Cells.Item(6, 5).Value = 19
$ws.Cells.Item(6, 6).Value = 0.07109340416113324
$ws.Cells.Item(6, 7).Value = 23

$ws.Cells.Item(7, 1).Value = 'Bayer 04 Leverkusen'
$ws.Cells.Item(7, 2).Value = 5.95253164556962
$ws.Cells.Item(7, 3).Value = 10.66846846846847
$ws.Cells.Item(7, 4).Value = 0.6293548387096775
$ws.Cells.Item(7, 5).Value = 18
$ws.Cells.Item(7, 6).Value = 0.05695773374519702
$ws.Cells.Item(7, 7).Value = 30

$ws.Cells.Item(8, 1).Value = 'Eintracht Frankfurt'
$ws.Cells.Item(8, 2).Value = 5.032846715328467
$ws.Cells.Item(8, 3).Value = 8.712707182320441
$ws.Cells.Item(8, 4).Value = 0.5190810946522721
$ws.Cells.Item(8, 5).Value = 21
$ws.Cells.Item(8, 6).Value = 0.07647058823529412
$ws.Cells.Item(8, 7).Value = -3

$ws.Cells.Item(9, 1).Value = 'SC Freiburg'
$ws.Cells.Item(9, 2).Value = 5.659634317862166
$ws.Cells.Item(9, 3).Value = 7.827586206896552
$ws.Cells.Item(9, 4).Value = 0.4634441901171402
$ws.Cells.Item(9, 5).Value = 26
$ws.Cells.Item(9, 6).Value = 0.07758811793393926
$ws.Cells.Item(9, 7).Value = -4

$ws.Cells.Item(10, 1).Value = 'Hamburger SV'
$ws.Cells.Item(10, 2).Value = 4.889795918367347
$ws.Cells.Item(10, 3).Value = 9.239406779661017
$ws.Cells.Item(10, 4).Value = 0.3734809613826627
$ws.Cells.Item(10, 5).Value = 15
$ws.Cells.Item(10, 6).Value = 0.08522400092134055
$ws.Cells.Item(10, 7).Value = -17

$ws.Cells.Item(11, 1).Value = '1. FC Union Berlin'
$ws.Cells.Item(11, 2).Value = 6.347665847665848
$ws.Cells.Item(11, 3).Value = 5.890109890109891
$ws.Cells.Item(11, 4).Value = 0.3512694748990191
$ws.Cells.Item(11, 5).Value = 9
$ws.Cells.Item(11, 6).Value = 0.13529081202585
$ws.Cells.Item(11, 7).Value = -1

$ws.Cells.Item(12, 1).Value = 'FC Augsburg'
$ws.Cells.Item(12, 2).Value = 6.139064475347661
$ws.Cells.Item(12, 3).Value = 6.435897435897436
$ws.Cells.Item(12, 4).Value = 0.4537301793605407
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 0.08169755455642067
$ws.Cells.Item(12, 7).Value = -25

$ws.Cells.Item(13, 1).Value = '1. FC Köln'
$ws.Cells.Item(13, 2).Value = 6.418154761904762
$ws.Cells.Item(13, 3).Value = 8.719298245614034
$ws.Cells.Item(13, 4).Value = 0.4478956568508807
$ws.Cells.Item(13, 5).Value = 12
$ws.Cells.Item(13, 6).Value = 0.07765486725663717
$ws.Cells.Item(13, 7).Value = -7

$ws.Cells.Item(14, 1).Value = 'Borussia M''gladbach'
$ws.Cells.Item(14, 2).Value = 5.821301775147929
$ws.Cells.Item(14, 3).Value = 9.145129224652088
$ws.Cells.Item(14, 4).Value = 0.3837955535182214
$ws.Cells.Item(14, 5).Value = 13
$ws.Cells.Item(14, 6).Value = 0.07205240174672489
$ws.Cells.Item(14, 7).Value = -11

$ws.Cells.Item(15, 1).Value = '1. FSV Mainz 05'
$ws.Cells.Item(15, 2).Value = 5.252475247524752
$ws.Cells.Item(15, 3).Value = 5.6875
$ws.Cells.Item(15, 4).Value = 0.4228187919463087
$ws.Cells.Item(15, 5).Value = 12
$ws.Cells.Item(15, 6).Value = 0.1080426356589147
$ws.Cells.Item(15, 7).Value = -5

$ws.Cells.Item(16, 1).Value = 'VfL Wolfsburg'
$ws.Cells.Item(16, 2).Value = 5.378002528445006
$ws.Cells.Item(16, 3).Value = 8.91194968553459
$ws.Cells.Item(16, 4).Value = 0.4011269820469139
$ws.Cells.Item(16, 5).Value = 12
$ws.Cells.Item(16, 6).Value = 0.09280717096633144
$ws.Cells.Item(16, 7).Value = -20

$ws.Cells.Item(17, 1).Value = 'SV Werder Bremen'
$ws.Cells.Item(17, 2).Value = 5.15158924205379
$ws.Cells.Item(17, 3).Value = 8.039182282793867
$ws.Cells.Item(17, 4).Value = 0.4706162832280367
$ws.Cells.Item(17, 5).Value = 17
$ws.Cells.Item(17, 6).Value = 0.09264356905552919
$ws.Cells.Item(17, 7).Value = -29

$ws.Cells.Item(18, 1).Value = 'FC St. Pauli'
$ws.Cells.Item(18, 2).Value = 6.943488943488943
$ws.Cells.Item(18, 3).Value = 9.004395604395604
$ws.Cells.Item(18, 4).Value = 0.4341364017739422
$ws.Cells.Item(18, 5).Value = 14
$ws.Cells.Item(18, 6).Value = 0.09461530601818739
$ws.Cells.Item(18, 7).Value = -33

$ws.Cells.Item(19, 1).Value = '1. FC Heidenheim'
$ws.Cells.Item(19, 2).Value = 5.647905759162303
$ws.Cells.Item(19, 3).Value = 8.47877358490566
$ws.Cells.Item(19, 4).Value = 0.3294277257030094
$ws.Cells.Item(19, 5).Value = 20
$ws.Cells.Item(19, 6).Value = 0.09476278496611214
$ws.Cells.Item(19, 7).Value = -18
